$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.964.79'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.288.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.47'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.15%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.857.66'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.65%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.59'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '68.899.88'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.254.59'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.87'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.25%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '395.87'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.35%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.76'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.522'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.61%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.75'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.11%  '

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.99'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'NEARProtocol'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.75'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.11%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.71%  '

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.31'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.60%  '

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.20'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.92%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.05%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.07'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.98'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.09%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.830'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.90'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.63'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.65'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.81%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.09%  '

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '25.67'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.36%  '

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.43'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0693'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.656.37'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '343.16'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0284'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '32.27'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.38'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.42%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.20%  '
